$d = $word.ActiveDocument

function Get-ParaRange($doc, $index) {
    $p = $doc.Paragraphs.Item($index)
    return $doc.Range($p.Range.Start, $p.Range.End)
}

# 1. Replace the director's full name in the preamble paragraph (paragraph 5).
$r1 = Get-ParaRange $d 5
$r1.Find.Execute("Kamolova Hulkar Ismoilovna", $true, $false, $false, $false, $false,
                  $true, 0, $false, "Murodova Aziza Baxtiyorovna", 2)

# 2. Merge the leading-space run with the "Professional malaka ..." run
#    on the "Nomi:" line (paragraph 87) into a single run.
$r2 = Get-ParaRange $d 87
$r2.Find.Execute(" Professional malaka oshirish nodavlat taʼlim muassasasi", $true, $false, $false, $false, $false,
                  $true, 0, $false, " Professional malaka oshirish nodavlat taʼlim muassasasi", 2)

# 3. Replace the signature initials/surname (paragraph 96).
$r3 = Get-ParaRange $d 96
$r3.Find.Execute("H.Kamolova", $true, $false, $false, $false, $false,
                  $true, 0, $false, "A. Murodova", 2)
